$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad / changed date) for rows 2-6 from 45184 to 45185
$ws.Range("C2").Value = 45185
$ws.Range("C3").Value = 45185
$ws.Range("C4").Value = 45185
$ws.Range("C5").Value = 45185
$ws.Range("C6").Value = 45185
